$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in missing A column values for rows 21-30
for ($r = 21; $r -le 30; $r++) {
    $ws.Cells.Item($r, 1).Value = $r
}

# Add summary formulas in D8/E8
$ws.Range("D8").Formula = "=SUM(B8,B9,B10,B11,B12,B13,B14,B15,B16,B17,B18,B19,B20,B21,B22,B23,B24,B25,B26,B27,B28,B29,B30)"
$ws.Range("E8").Formula = "=(D8/25)"

# Re-enter B100's formula explicitly so it is its own (non-shared) formula
# rather than being auto-filled as part of the B67:B100 shared formula group.
$ws.Range("B100").Formula = "=PRODUCT(B99,C99)"

# Update the view: selection on E5 (also resets the scrolled top-left cell)
$ws.Range("E5").Select()
